$d = $word.ActiveDocument

# ============================================================
# 1. Underline the section title "7. Recommendations for Future
#    Initiatives for Japan to improve Digital Transformation in
#    the Electronics/Semiconductor industry."
# ============================================================
$d.Paragraphs.Item(1).Range.Font.Underline = 1

# ============================================================
# 2. Collapse the 3 runs that spell out the "5 key factors (...)"
#    sentence into a single run (pure re-flow, text unchanged).
# ============================================================
$old2 = "Technology, Organization Design, People, Leadership, Culture"
$rng2 = $d.Content
$rng2.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $old2, 2) | Out-Null

# ============================================================
# 3. Mark the run that carries the inline picture as NoProof
#    (adds <w:noProof/> to its run properties).
# ============================================================
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ip = $d.Paragraphs.Item($i)
    if ($ip.Range.InlineShapes.Count -gt 0) {
        $ip.Range.InlineShapes.Item(1).Range.NoProofing = 1
    }
}

# ============================================================
# 4. Split the LSTC sentence to add:
#    ", or any other R&D centers within Japanese semiconductors"
#    right after "LSTC (being a R&D center)".
# ============================================================
$rng4 = $d.Content
$rng4.Find.Execute("LSTC (being a R&D center)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng4.Collapse(0)
$rng4.InsertAfter(", or any other R&D centers within Japanese semiconductors")
$rng4.Bold = 1
$rng4.Bold = 0

# ============================================================
# 5. Split the sentence right before "for mass production" to add
#    "or any other HVM manufacturing ".
# ============================================================
$rng5 = $d.Content
$rng5.Find.Execute(" for mass production/commercialization. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPoint5 = $d.Range($rng5.Start + 1, $rng5.Start + 1)
$insertPoint5.InsertAfter("or any other HVM manufacturing ")
$insertPoint5.Bold = 1
$insertPoint5.Bold = 0

# ============================================================
# 6. Append six new paragraphs (Organization Design / People /
#    Leadership / Culture dimension sections) right before the
#    trailing empty paragraph at the end of the document.
# ============================================================
function New-TailParagraph {
    $lastIdx = $d.Paragraphs.Count
    $a = $d.Paragraphs.Item($lastIdx)
    $insPt = $d.Range($a.Range.Start, $a.Range.Start)
    $insPt.InsertParagraphBefore()
    $newIdx = $d.Paragraphs.Count - 1
    return $d.Paragraphs.Item($newIdx)
}

# --- paragraph 1 ---
$p = New-TailParagraph
$p.Range.Text = "Organization Design Dimension of Agility:"
$p.Range.Font.Underline = 1

# --- paragraph 2 ---
$p = New-TailParagraph
$p.Range.Text = "In most cases, Japanese electronics/semiconductor firms’ organization structure are highly functional structure and this leads to organization silos within the company. Let us take the case of a Japanese electronics firm JE1 firm (source). JE1 was founded in the 1940s, and initially the firm has successively introduced innovative finished electronic products into the global market. However, the performance of the JE1 firm start to declined after the early 2000s. This is largely due to the outsourcing of core component production after the early 2000s, where JE1 decide to close down their manufacturing engineering divisions and create a new production subsidiary. JE1 still retain the R&D divisions and product design engineering divisions. The production will be done by a production subsidiary that includes production process engineering and production divisions. This lead to organization silos and created a large barriers between product design engineering divisions at the JE1 firm and its production subsidiary process engineering sections. Prior to the establishment of the production subsidiary, all employees across the product design and production divisions share common corporate values and company mission values. After the establishment of the production subsidiary, the product design engineers in JE1 become much less involved in the production process and interacted less with the production process engineers. Vice versa, production workers and supervisors at the production subsidiary do not propose many ideas based on their work experience to the product design engineers as before. This organization silos led to the decline of technological superiority. From this case study, it is quite clear that JE1 organization structure is very functional. Employees within each functional area (product design, production etc.) has a very deep level of knowledge within their domain. However, the drawback is that it creates distinct silos and makes collaboration between different specialty areas difficult as evidenced by the JE1 firm. Having this clear delineation and separation has a tendency to slow down decision making and makes adapting to market conditions challenging (therefore JE1 lost the technological innovation superiority that they had held previously)."
$rng = $p.Range.Duplicate
$rng.Find.Execute("source", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Font.Italic = 1
$rng.Font.ItalicBi = 1

# --- paragraph 3 ---
$p = New-TailParagraph
$p.Range.Text = "To avoid such organizational design pitfalls that JE1 faced, the future of Japanese electronics should adopt a matrix structure. A matrix structure can be seen as a hybrid combination of functional and divisional structures. Unlike functional structure where it is optimizes for a given function and a divisional structure that optimizes for a particular product/service, a matrix structure aims to utilizes a company’s resources and assets and people towards a shared common goal. The key benefits of a matrix organization is that it drives an increased level of collaboration and communication across the organization. People are not bounded to a given function and instead, people work across the boundaries to achieve the same goal regardless of where they are in the organization chart. In the case of JE1, if the firm had earlier on adopt a matrix structure which encourages collaboration between product design engineers and production engineers, it might lead to more innovative processes and eventually lead to products for them to remain technologically innovative. However, that being said, a matrix structure can also have its drawbacks whereby people have to report to multiple project managers at the same time. One of the solution to tackle this is the prioritization of projects e.g. which project has the most impact on business values and therefore people are able to collaborate closely and achieve the best business value outcome. "
$splitPt = $p.Range.Duplicate
$splitPt.Find.Execute("To avoid such organizational design pitfalls that JE1 faced, the future of Japanese electronics should adopt a matrix structure. A matrix structure can be seen as a hybrid combination of functional and ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPt.Collapse(0)
$splitPt.Bold = 1
$splitPt.Bold = 0
$splitPt = $p.Range.Duplicate
$splitPt.Find.Execute("To avoid such organizational design pitfalls that JE1 faced, the future of Japanese electronics should adopt a matrix structure. A matrix structure can be seen as a hybrid combination of functional and divisional structures. Unlike functional structure where it is optimizes for a given function and a divisional structure that optimizes for a particular product/service, a matrix structure aims to utilizes a company’s resources and assets and people towards a shared common goal. ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPt.Collapse(0)
$splitPt.Bold = 1
$splitPt.Bold = 0

# --- paragraph 4 ---
$p = New-TailParagraph
$p.Range.Text = "People Dimension of Agility:"
$p.Range.Font.Underline = 1

# --- paragraph 5 ---
$p = New-TailParagraph
$p.Range.Text = "Leadership Dimension of Agility:"
$p.Range.Font.Underline = 1

# --- paragraph 6 ---
$p = New-TailParagraph
$p.Range.Text = "Culture Dimension of Agility:"
$p.Range.Font.Underline = 1
